$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 376.6154
$ws.Range("I6").Value = 144
$ws.Range("J6").Value = 900
$ws.Range("K6").Value = 432
$ws.Range("L6").Value = 2700
$ws.Range("M6").Value = -320
$ws.Range("N6").Value = -2924
$ws.Range("H31").Value = 1545.6
$ws.Range("I31").Value = 1500.75
$ws.Range("J31").Value = 1725
$ws.Range("K31").Value = 4502.25
$ws.Range("L31").Value = 5175
$ws.Range("M31").Value = -4272.25
$ws.Range("N31").Value = -5635
$ws.Range("H112").Value = 1123.8182
$ws.Range("J112").Value = 1144
$ws.Range("L112").Value = 3432
$ws.Range("N112").Value = -5648
$ws.Range("H137").Value = 1889.95
$ws.Range("I137").Value = 1190.909
$ws.Range("K137").Value = 3572.727
$ws.Range("M137").Value = -1022.727
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23909.705
$ws.Range("I32").Value = 4415.875
$ws.Range("K32").Value = 4415.875
$ws.Range("M32").Value = -4128.875
$ws.Range("H44").Value = 14985.714
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 14985.714
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 14985.714
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -15961.714
$ws.Range("H61").Value = 1460.8833
$ws.Range("I61").Value = 889.4857
$ws.Range("K61").Value = 889.4857
$ws.Range("M61").Value = -677.4857
$ws.Range("H63").Value = 1661.2
$ws.Range("I63").Value = 1300
$ws.Range("J63").Value = 2203
$ws.Range("K63").Value = 1300
$ws.Range("L63").Value = 2203
$ws.Range("M63").Value = -614
$ws.Range("N63").Value = -3575
$ws.Range("H66").Value = 1661.2
$ws.Range("I66").Value = 1300
$ws.Range("J66").Value = 2203
$ws.Range("K66").Value = 6500
$ws.Range("L66").Value = 11015
$ws.Range("M66").Value = -3068
$ws.Range("N66").Value = -17879
$ws.Range("H74").Value = 753.2143
$ws.Range("I74").Value = 705
$ws.Range("K74").Value = 705
$ws.Range("M74").Value = 169
$ws.Range("H77").Value = 753.2143
$ws.Range("I77").Value = 705
$ws.Range("K77").Value = 3525
$ws.Range("M77").Value = 843
$ws.Range("H136").Value = 1460.8833
$ws.Range("I136").Value = 889.4857
$ws.Range("K136").Value = 2668.4571
$ws.Range("M136").Value = -118.4570999999996
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3800.2856
$ws.Range("I134").Value = 3240.3
$ws.Range("K134").Value = 9720.900000000001
$ws.Range("M134").Value = -7185.900000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38109.676
$ws.Range("I31").Value = 2249.5
$ws.Range("K31").Value = 2249.5
$ws.Range("M31").Value = -1954.5
$ws.Range("H34").Value = 38109.676
$ws.Range("I34").Value = 2249.5
$ws.Range("K34").Value = 2249.5
$ws.Range("M34").Value = -2047.5
$ws.Range("H58").Value = 1423.2368
$ws.Range("I58").Value = 1275
$ws.Range("J58").Value = 2079.7144
$ws.Range("K58").Value = 1275
$ws.Range("L58").Value = 2079.7144
$ws.Range("M58").Value = -1072
$ws.Range("N58").Value = -2485.7144
$ws.Range("H64").Value = 34660
$ws.Range("J64").Value = 34660
$ws.Range("L64").Value = 34660
$ws.Range("N64").Value = -35156
$ws.Range("H67").Value = 34660
$ws.Range("J67").Value = 34660
$ws.Range("L67").Value = 34660
$ws.Range("N67").Value = -36376
$ws.Range("H107").Value = 6779
$ws.Range("I107").Value = 13454.875
$ws.Range("J107").Value = 844.8889
$ws.Range("K107").Value = 13454.875
$ws.Range("L107").Value = 844.8889
$ws.Range("M107").Value = -11534.875
$ws.Range("N107").Value = -4684.8889
$ws.Range("H131").Value = 18453.055
$ws.Range("J131").Value = 18453.055
$ws.Range("L131").Value = 18453.055
$ws.Range("N131").Value = -28533.055
$ws.Range("H132").Value = 3459.0967
$ws.Range("I132").Value = 3506.6
$ws.Range("J132").Value = 3372.7273
$ws.Range("K132").Value = 10519.8
$ws.Range("L132").Value = 10118.1819
$ws.Range("M132").Value = -7989.799999999999
$ws.Range("N132").Value = -15178.1819
$ws.Range("H136").Value = 1423.2368
$ws.Range("I136").Value = 1275
$ws.Range("J136").Value = 2079.7144
$ws.Range("K136").Value = 3825
$ws.Range("L136").Value = 6239.1432
$ws.Range("M136").Value = -1275
$ws.Range("N136").Value = -11339.1432
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1750
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 1500
$ws.Range("K32").Value = 6000
$ws.Range("L32").Value = 4500
$ws.Range("M32").Value = -5717
$ws.Range("N32").Value = -5066
$ws.Range("H39").Value = 1950
$ws.Range("J39").Value = 2916.6667
$ws.Range("L39").Value = 8750.000100000001
$ws.Range("N39").Value = -9338.000100000001
$ws.Range("H55").Value = 12069.857
$ws.Range("I55").Value = 33733.332
$ws.Range("J55").Value = 8459.277
$ws.Range("K55").Value = 101199.996
$ws.Range("L55").Value = 25377.831
$ws.Range("M55").Value = -101022.996
$ws.Range("N55").Value = -25731.831
$ws.Range("H127").Value = 861.25
$ws.Range("J127").Value = 861.25
$ws.Range("L127").Value = 2583.75
$ws.Range("N127").Value = -12503.75
$ws.Range("H140").Value = 1330.8438
$ws.Range("I140").Value = 944.3333
$ws.Range("J140").Value = 2068.7273
$ws.Range("K140").Value = 2832.9999
$ws.Range("L140").Value = 6206.1819
$ws.Range("M140").Value = 2347.0001
$ws.Range("N140").Value = -16566.1819
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 200
$ws.Range("J4").Value = 200
$ws.Range("L4").Value = 200
$ws.Range("N4").Value = -424
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 20000
$ws.Range("J5").Value = 20000
$ws.Range("L5").Value = 20000
$ws.Range("N5").Value = -20226
$ws.Range("H61").Value = 1385.64
$ws.Range("I61").Value = 1215.4445
$ws.Range("K61").Value = 1215.4445
$ws.Range("M61").Value = -1013.4445
$ws.Range("H93").Value = 3779
$ws.Range("I93").Value = 3722.75
$ws.Range("K93").Value = 3722.75
$ws.Range("M93").Value = -2474.75
$ws.Range("H113").Value = 1385.64
$ws.Range("I113").Value = 1215.4445
$ws.Range("K113").Value = 1215.4445
$ws.Range("M113").Value = 954.5554999999999
$ws.Range("H122").Value = 2448.7856
$ws.Range("I122").Value = 2478.2222
$ws.Range("J122").Value = 2395.8
$ws.Range("K122").Value = 7434.6666
$ws.Range("L122").Value = 7187.400000000001
$ws.Range("M122").Value = -4984.6666
$ws.Range("N122").Value = -12087.4
$ws.Range("H132").Value = 1992.804
$ws.Range("I132").Value = 2141.0435
$ws.Range("J132").Value = 629
$ws.Range("K132").Value = 6423.130500000001
$ws.Range("L132").Value = 1887
$ws.Range("M132").Value = -3893.130500000001
$ws.Range("N132").Value = -6947
$ws.Range("H136").Value = 967.3333
$ws.Range("I136").Value = 926.35895
$ws.Range("K136").Value = 2779.07685
$ws.Range("M136").Value = -229.0768500000004
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 7500
$ws.Range("J19").Value = 7500
$ws.Range("L19").Value = 7500
$ws.Range("N19").Value = -7848
$ws.Range("H100").Value = 63237.688
$ws.Range("I100").Value = 125550.5
$ws.Range("J100").Value = 924.875
$ws.Range("K100").Value = 251101
$ws.Range("L100").Value = 1849.75
$ws.Range("M100").Value = -250560
$ws.Range("N100").Value = -2931.75
$ws.Range("H132").Value = 1664.2344
$ws.Range("I132").Value = 1677.5593
$ws.Range("J132").Value = 1507
$ws.Range("K132").Value = 5032.6779
$ws.Range("L132").Value = 4521
$ws.Range("M132").Value = -2502.6779
$ws.Range("N132").Value = -9581
$ws.Range("H136").Value = 591.9056399999999
$ws.Range("I136").Value = 354.51282
$ws.Range("J136").Value = 1253.2142
$ws.Range("K136").Value = 1063.53846
$ws.Range("L136").Value = 3759.6426
$ws.Range("M136").Value = 1486.46154
$ws.Range("N136").Value = -8859.642599999999
